# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded on the handback-status
# report: the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Correspond Handoff/Handback Datetime" columns on the per-locale sheets
# for the b1a02e00-... file.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 05:10:21"

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 05:10:17"
$wsZhCn.Range("K2").Value = "2016-09-01 05:10:34"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 05:10:21"
$wsDeDe.Range("K2").Value = "2016-09-01 05:10:41"
